$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of data appended at row 92, mirroring the existing rows'
# "temperature" feed readings pulled from Adafruit IO.
$row = 92

# Force text formatting first so numeric-looking strings (e.g. "25")
# are preserved as text rather than being coerced to numbers.
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 2).NumberFormat = "@"
$ws.Cells.Item($row, 3).NumberFormat = "@"
$ws.Cells.Item($row, 4).NumberFormat = "@"
$ws.Cells.Item($row, 5).NumberFormat = "@"
$ws.Cells.Item($row, 6).NumberFormat = "@"

$ws.Cells.Item($row, 1).Value = "2024-09-25T18:06:40Z"
$ws.Cells.Item($row, 2).Value = "temperature"
$ws.Cells.Item($row, 3).Value = "25"
$ws.Cells.Item($row, 4).Value = "N/A"
$ws.Cells.Item($row, 5).Value = "N/A"
$ws.Cells.Item($row, 6).Value = "N/A"
